$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 116735
$ws.Range("J17").Value = 116735
$ws.Range("L17").Value = 350205
$ws.Range("N17").Value = -350541
$ws.Range("H53").Value = 972.7
$ws.Range("I53").Value = 1114.9412
$ws.Range("J53").Value = 166.66667
$ws.Range("K53").Value = 1114.9412
$ws.Range("L53").Value = 166.66667
$ws.Range("M53").Value = -477.9412
$ws.Range("N53").Value = -1440.66667
$ws.Range("H107").Value = 614.4167
$ws.Range("I107").Value = 639.3
$ws.Range("J107").Value = 490
$ws.Range("K107").Value = 639.3
$ws.Range("L107").Value = 490
$ws.Range("M107").Value = 1280.7
$ws.Range("N107").Value = -4330
$ws.Range("H138").Value = 1737.55
$ws.Range("I138").Value = 735.6
$ws.Range("J138").Value = 2557.3274
$ws.Range("K138").Value = 2206.8
$ws.Range("L138").Value = 7671.9822
$ws.Range("M138").Value = 2933.2
$ws.Range("N138").Value = -17951.9822

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4734.5137
$ws.Range("I32").Value = 2934.869
$ws.Range("J32").Value = 14714.363
$ws.Range("K32").Value = 2934.869
$ws.Range("L32").Value = 14714.363
$ws.Range("M32").Value = -2647.869
$ws.Range("N32").Value = -15288.363
$ws.Range("H74").Value = 68450.39999999999
$ws.Range("I74").Value = 85004.664
$ws.Range("J74").Value = 2233.3333
$ws.Range("K74").Value = 85004.664
$ws.Range("L74").Value = 2233.3333
$ws.Range("M74").Value = -84130.664
$ws.Range("N74").Value = -3981.3333
$ws.Range("H77").Value = 68450.39999999999
$ws.Range("I77").Value = 85004.664
$ws.Range("J77").Value = 2233.3333
$ws.Range("K77").Value = 425023.32
$ws.Range("L77").Value = 11166.6665
$ws.Range("M77").Value = -420655.32
$ws.Range("N77").Value = -19902.6665
$ws.Range("H110").Value = 1478.7333
$ws.Range("I110").Value = 1598.5385
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 1598.5385
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 446.4614999999999
$ws.Range("N110").Value = -4790
$ws.Range("H133").Value = 31999.666
$ws.Range("J133").Value = 31999.666
$ws.Range("L133").Value = 31999.666
$ws.Range("N133").Value = -37059.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5913.057
$ws.Range("I134").Value = 5711.619
$ws.Range("J134").Value = 6215.2144
$ws.Range("K134").Value = 17134.857
$ws.Range("L134").Value = 18645.6432
$ws.Range("M134").Value = -14599.857
$ws.Range("N134").Value = -23715.6432
$ws.Range("H137").Value = 39409
$ws.Range("J137").Value = 39409
$ws.Range("L137").Value = 39409
$ws.Range("N137").Value = -49609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35716024
$ws.Range("I31").Value = 55556668
$ws.Range("J31").Value = 2866.8
$ws.Range("K31").Value = 55556668
$ws.Range("L31").Value = 2866.8
$ws.Range("M31").Value = -55556373
$ws.Range("N31").Value = -3456.8
$ws.Range("H34").Value = 35716024
$ws.Range("I34").Value = 55556668
$ws.Range("J34").Value = 2866.8
$ws.Range("K34").Value = 55556668
$ws.Range("L34").Value = 2866.8
$ws.Range("M34").Value = -55556466
$ws.Range("N34").Value = -3270.8
$ws.Range("H107").Value = 1675.3914
$ws.Range("I107").Value = 448.7857
$ws.Range("J107").Value = 3583.4443
$ws.Range("K107").Value = 448.7857
$ws.Range("L107").Value = 3583.4443
$ws.Range("M107").Value = 1471.2143
$ws.Range("N107").Value = -7423.4443
$ws.Range("H132").Value = 1847.8948
$ws.Range("I132").Value = 1245.2094
$ws.Range("J132").Value = 3699
$ws.Range("K132").Value = 3735.6282
$ws.Range("L132").Value = 11097
$ws.Range("M132").Value = -1205.6282
$ws.Range("N132").Value = -16157
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
$ws.Range("H137").Value = 38950
$ws.Range("J137").Value = 38950
$ws.Range("L137").Value = 38950
$ws.Range("N137").Value = -49150

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4180
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 4975
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 14925
$ws.Range("M22").Value = -2831
$ws.Range("N22").Value = -15263
$ws.Range("H27").Value = 4180
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 4975
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 14925
$ws.Range("M27").Value = -2898
$ws.Range("N27").Value = -15129
$ws.Range("H64").Value = 5265
$ws.Range("I64").Value = 3584.6667
$ws.Range("K64").Value = 10754.0001
$ws.Range("M64").Value = -10484.0001
$ws.Range("H67").Value = 5265
$ws.Range("I67").Value = 3584.6667
$ws.Range("K67").Value = 10754.0001
$ws.Range("M67").Value = -9818.000100000001
$ws.Range("H75").Value = 2634.923
$ws.Range("I75").Value = 1200
$ws.Range("J75").Value = 3531.75
$ws.Range("K75").Value = 3600
$ws.Range("L75").Value = 10595.25
$ws.Range("M75").Value = -2602
$ws.Range("N75").Value = -12591.25
$ws.Range("H78").Value = 2634.923
$ws.Range("I78").Value = 1200
$ws.Range("J78").Value = 3531.75
$ws.Range("K78").Value = 10800
$ws.Range("L78").Value = 31785.75
$ws.Range("M78").Value = -5808
$ws.Range("N78").Value = -41769.75
$ws.Range("H98").Value = 1599.3334
$ws.Range("I98").Value = 900
$ws.Range("J98").Value = 1949
$ws.Range("K98").Value = 2700
$ws.Range("L98").Value = 5847
$ws.Range("M98").Value = -1202
$ws.Range("N98").Value = -8843
$ws.Range("H105").Value = 166103.62
$ws.Range("J105").Value = 166103.62
$ws.Range("L105").Value = 498310.86
$ws.Range("N105").Value = -503552.86
$ws.Range("H113").Value = 1515682.2
$ws.Range("I113").Value = 6061151
$ws.Range("J113").Value = 526
$ws.Range("K113").Value = 18183453
$ws.Range("L113").Value = 1578
$ws.Range("M113").Value = -18181283
$ws.Range("N113").Value = -5918
$ws.Range("H114").Value = 7637237
$ws.Range("J114").Value = 5683006.5
$ws.Range("L114").Value = 17049019.5
$ws.Range("N114").Value = -17055527.5
$ws.Range("H120").Value = 8400
$ws.Range("I120").Value = 4800
$ws.Range("K120").Value = 14400
$ws.Range("M120").Value = -9562
$ws.Range("H129").Value = 2505.639
$ws.Range("I129").Value = 1700.2222
$ws.Range("J129").Value = 3311.0557
$ws.Range("K129").Value = 5100.6666
$ws.Range("L129").Value = 9933.167099999999
$ws.Range("M129").Value = -100.6665999999996
$ws.Range("N129").Value = -19933.1671
$ws.Range("H131").Value = 864.53
$ws.Range("J131").Value = 909.375
$ws.Range("L131").Value = 2728.125
$ws.Range("N131").Value = -12808.125
$ws.Range("H137").Value = 30763140
$ws.Range("I137").Value = 3048.75
$ws.Range("J137").Value = 48340336
$ws.Range("K137").Value = 9146.25
$ws.Range("L137").Value = 145021008
$ws.Range("M137").Value = -4046.25
$ws.Range("N137").Value = -145031208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1673002.9
$ws.Range("I14").Value = 2006602.4
$ws.Range("J14").Value = 5005
$ws.Range("K14").Value = 2006602.4
$ws.Range("L14").Value = 5005
$ws.Range("M14").Value = -2006434.4
$ws.Range("N14").Value = -5341
$ws.Range("H80").Value = 2589.7368
$ws.Range("I80").Value = 2569.0625
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 2569.0625
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -1571.0625
$ws.Range("N80").Value = -4696
$ws.Range("H83").Value = 2589.7368
$ws.Range("I83").Value = 2569.0625
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 12845.3125
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -7853.3125
$ws.Range("N83").Value = -23484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 607.1177
$ws.Range("I16").Value = 582.5625
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 582.5625
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -412.5625
$ws.Range("N16").Value = -1340
$ws.Range("H22").Value = 553.9394
$ws.Range("I22").Value = 303.33334
$ws.Range("J22").Value = 609.62964
$ws.Range("K22").Value = 303.33334
$ws.Range("L22").Value = 609.62964
$ws.Range("M22").Value = -8.333340000000021
$ws.Range("N22").Value = -1199.62964
$ws.Range("H27").Value = 553.9394
$ws.Range("I27").Value = 303.33334
$ws.Range("J27").Value = 609.62964
$ws.Range("K27").Value = 303.33334
$ws.Range("L27").Value = 609.62964
$ws.Range("M27").Value = -196.33334
$ws.Range("N27").Value = -823.62964
$ws.Range("H46").Value = 920.53845
$ws.Range("I46").Value = 1433.6666
$ws.Range("J46").Value = 766.6
$ws.Range("K46").Value = 1433.6666
$ws.Range("L46").Value = 766.6
$ws.Range("M46").Value = -1245.6666
$ws.Range("N46").Value = -1142.6
$ws.Range("H55").Value = 544.95
$ws.Range("I55").Value = 168.75
$ws.Range("K55").Value = 168.75
$ws.Range("M55").Value = 4.25
